$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.479.75'
$ws.Range("E2").Value = '  -2.23%  '
$ws.Range("D3").Value = '2.288.85'
$ws.Range("E3").Value = '  -1.07%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.00'
$ws.Range("E5").Value = '  -2.86%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.17'
$ws.Range("E6").Value = '  -6.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.504'
$ws.Range("E7").Value = '  -5.56%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  -5.92%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.25'
$ws.Range("E10").Value = '  -6.82%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0787'
$ws.Range("E11").Value = '  -3.23%  '
$ws.Range("E12").Value = '  -0.02%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.70'
$ws.Range("E13").Value = '  -4.46%  '
$ws.Range("D14").Value = '2.643.43'
$ws.Range("E14").Value = '  -1.13%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.59'
$ws.Range("E15").Value = '  +2.57%  '
$ws.Range("D16").Value = '2.282.78'
$ws.Range("E16").Value = '  -1.42%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.802'
$ws.Range("E17").Value = '  -1.17%  '
$ws.Range("D18").Value = '42.402.44'
$ws.Range("E18").Value = '  -2.19%  '
$ws.Range("D19").Value = '0.0₃0897'
$ws.Range("E19").Value = '  -3.10%  '
$ws.Range("B20").Value = 'InternetComputer(DFINITY)'
$ws.Range("C20").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.48'
$ws.Range("E20").Value = '  -5.66%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.02'
$ws.Range("E21").Value = '  -2.37%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.83'
$ws.Range("E22").Value = '  -0.64%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '234.41'
$ws.Range("E23").Value = '  -3.30%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.97'
$ws.Range("E24").Value = '  -3.03%  '
$ws.Range("E25").Value = '  -4.04%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.96'
$ws.Range("E27").Value = '  +0.44%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.17'
$ws.Range("E28").Value = '  -9.31%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '34.62'
$ws.Range("E29").Value = '  -7.21%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.14'
$ws.Range("E30").Value = '  -5.19%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '162.09'
$ws.Range("E31").Value = '  -3.07%  '
$ws.Range("E32").Value = '  -0.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.01'
$ws.Range("E33").Value = '  -5.46%  '
$ws.Range("E34").Value = '  +1.90%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.41'
$ws.Range("E35").Value = '  -5.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0712'
$ws.Range("E36").Value = '  -4.43%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '16.96'
$ws.Range("E37").Value = '  -7.51%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.88'
$ws.Range("E38").Value = '  -6.01%  '
$ws.Range("E39").Value = '  -4.80%  '
$ws.Range("E40").Value = '  -5.40%  '
$ws.Range("E41").Value = '  -4.14%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.45'
$ws.Range("E42").Value = '  -9.37%  '
$ws.Range("D43").Value = '1.972.36'
$ws.Range("E43").Value = '  -1.18%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '18.66'
$ws.Range("E44").Value = '  -2.26%  '
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0279'
$ws.Range("E45").Value = '  -4.99%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.09'
$ws.Range("E46").Value = '  +0.77%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.86'
$ws.Range("E47").Value = '  -7.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.27'
$ws.Range("E48").Value = '  -3.11%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.84'
$ws.Range("E49").Value = '  -3.65%  '
$ws.Range("D50").Value = '2.514.31'
$ws.Range("E50").Value = '  -0.99%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.65'
$ws.Range("E51").Value = '  -1.08%  '
